# The <id>…</id> run for this folio was split across three runs:
#   "<id>"  (Courier New / 7f6000 / 18)  +  "p085v_a4" (plain)  +  "</id>" (Courier New / 7f6000 / 18)
# The edit fixes the typo'd id (drop the stray "a") and collapses the three
# runs into a single run, inheriting the Courier-New/gold formatting that
# already bookended the text.
$d = $word.ActiveDocument

$old = "<id>p085v_a4</id>"
$new = "<id>p085v_4</id>"

$range = $d.Content
$found = $range.Find.Execute($old, $false, $false, $false, $false, $false, `
                              $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target text '$old' to replace."
}
